# ---------------------------------------------------------------------------
# Scheduled Universalis market-data refresh for Zodiark_Profits.xlsx
#
# Re-pulls current average NQ/HQ marketboard prices for the leve turn-in items
# on each job sheet and recomputes the dependent Leve price / profit columns
# (H:N) for every row whose source price moved since the previous run. Rows
# whose NQ or HQ price dropped out of the marketboard entirely have their now-
# meaningless profit cell cleared instead of zeroed.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===== ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 2695.88
$ws.Range("I98").Value = 1699.909
$ws.Range("J98").Value = 9999.666999999999
$ws.Range("K98").Value = 1699.909
$ws.Range("L98").Value = 9999.666999999999
$ws.Range("M98").Value = -201.9090000000001
$ws.Range("N98").Value = -12995.667
# Row 101
$ws.Range("H101").Value = 1183
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()
# Row 113
$ws.Range("H113").Value = 22046.072
$ws.Range("J113").Value = 8870.799999999999
$ws.Range("L113").Value = 8870.799999999999
$ws.Range("N113").Value = -15378.8
# Row 115
$ws.Range("H115").Value = 4950
$ws.Range("I115").Value = 4950
$ws.Range("K115").Value = 14850
$ws.Range("M115").Value = -13283
# Row 121
$ws.Range("H121").Value = 1733.7
$ws.Range("J121").Value = 1733.7
$ws.Range("L121").Value = 5201.1
$ws.Range("N121").Value = -8695.1
# Row 122
$ws.Range("H122").Value = 2695.88
$ws.Range("I122").Value = 1699.909
$ws.Range("J122").Value = 9999.666999999999
$ws.Range("K122").Value = 5099.727000000001
$ws.Range("L122").Value = 29999.001
$ws.Range("M122").Value = -2649.727000000001
$ws.Range("N122").Value = -34899.001
# Row 135
$ws.Range("H135").Value = 3106.0908
$ws.Range("I135").Value = 3106.0908
$ws.Range("K135").Value = 27954.8172
$ws.Range("M135").Value = -25419.8172
# Row 137
$ws.Range("H137").Value = 1115.1666
$ws.Range("J137").Value = 1172.75
$ws.Range("L137").Value = 3518.25
$ws.Range("N137").Value = -8618.25
# Row 138
$ws.Range("H138").Value = 4677.8057
$ws.Range("J138").Value = 6435.6
$ws.Range("L138").Value = 19306.8
$ws.Range("N138").Value = -29586.8

# ===== ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 1279.6
$ws.Range("I97").Value = 950
$ws.Range("J97").Value = 1499.3334
$ws.Range("K97").Value = 950
$ws.Range("L97").Value = 1499.3334
$ws.Range("M97").Value = -454
$ws.Range("N97").Value = -2491.3334
# Row 98
$ws.Range("H98").Value = 49614.75
$ws.Range("J98").Value = 49614.75
$ws.Range("L98").Value = 49614.75
$ws.Range("N98").Value = -55604.75
# Row 132
$ws.Range("H132").Value = 7463.0303
$ws.Range("I132").Value = 7431.3896
$ws.Range("K132").Value = 22294.1688
$ws.Range("M132").Value = -19764.1688

# ===== BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 16
$ws.Range("H16").Value = 100000
$ws.Range("J16").Value = 100000
$ws.Range("L16").Value = 100000
$ws.Range("N16").Value = -100340
# Row 20
$ws.Range("H20").Value = 990.8889
$ws.Range("I20").Value = 981.2692
$ws.Range("J20").Value = 1015.9
$ws.Range("K20").Value = 981.2692
$ws.Range("L20").Value = 1015.9
$ws.Range("M20").Value = -734.2692
$ws.Range("N20").Value = -1509.9
# Row 23
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
# Row 102
$ws.Range("H102").Value = 15271
$ws.Range("I102").Value = 15271
$ws.Range("K102").Value = 15271
$ws.Range("M102").Value = -12026

# ===== CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 7500350
$ws.Range("I6").Value = 7500350
$ws.Range("K6").Value = 7500350
$ws.Range("M6").Value = -7500237
# Row 99
$ws.Range("H99").Value = 1243.4839
$ws.Range("I99").Value = 1259.7587
$ws.Range("J99").Value = 1007.5
$ws.Range("K99").Value = 1259.7587
$ws.Range("L99").Value = 1007.5
$ws.Range("M99").Value = 238.2412999999999
$ws.Range("N99").Value = -4003.5
# Row 126
$ws.Range("H126").Value = 1243.4839
$ws.Range("I126").Value = 1259.7587
$ws.Range("J126").Value = 1007.5
$ws.Range("K126").Value = 3779.2761
$ws.Range("L126").Value = 3022.5
$ws.Range("M126").Value = -1309.2761
$ws.Range("N126").Value = -7962.5
# Row 132
$ws.Range("H132").Value = 2011.7826
$ws.Range("I132").Value = 1830.8235
$ws.Range("J132").Value = 2524.5
$ws.Range("K132").Value = 5492.470499999999
$ws.Range("L132").Value = 7573.5
$ws.Range("M132").Value = -2962.470499999999
$ws.Range("N132").Value = -12633.5

# ===== CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 28682074
$ws.Range("I4").Value = 35848104
$ws.Range("K4").Value = 107544312
$ws.Range("M4").Value = -107544200
# Row 81
$ws.Range("H81").Value = 7950
$ws.Range("I81").Value = 10875
$ws.Range("K81").Value = 32625
$ws.Range("M81").Value = -31502
# Row 84
$ws.Range("H84").Value = 7950
$ws.Range("I84").Value = 10875
$ws.Range("K84").Value = 97875
$ws.Range("M84").Value = -92259

# ===== GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 5171.9
$ws.Range("I122").Value = 4713.875
$ws.Range("J122").Value = 7004
$ws.Range("K122").Value = 14141.625
$ws.Range("L122").Value = 21012
$ws.Range("M122").Value = -11691.625
$ws.Range("N122").Value = -25912
# Row 126
$ws.Range("H126").Value = 65006668
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 65006668
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 195020004
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -195024944
# Row 132
$ws.Range("H132").Value = 17386.285
$ws.Range("I132").Value = 20341.4
$ws.Range("J132").Value = 9998.5
$ws.Range("K132").Value = 61024.2
$ws.Range("L132").Value = 29995.5
$ws.Range("M132").Value = -58494.2
$ws.Range("N132").Value = -35055.5

# ===== LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 17861364
$ws.Range("I7").Value = 3966
$ws.Range("J7").Value = 71433560
$ws.Range("K7").Value = 3966
$ws.Range("L7").Value = 71433560
$ws.Range("M7").Value = -3854
$ws.Range("N7").Value = -71433784
# Row 16
$ws.Range("H16").Value = 1684.3334
$ws.Range("I16").Value = 921.9524
$ws.Range("K16").Value = 921.9524
$ws.Range("M16").Value = -751.9524
# Row 40
$ws.Range("H40").Value = 6310.0713
$ws.Range("I40").Value = 6294
$ws.Range("J40").Value = 6331.5
$ws.Range("K40").Value = 6294
$ws.Range("L40").Value = 6331.5
$ws.Range("M40").Value = -6158
$ws.Range("N40").Value = -6603.5
# Row 126
$ws.Range("H126").Value = 17861364
$ws.Range("I126").Value = 3966
$ws.Range("J126").Value = 71433560
$ws.Range("K126").Value = 11898
$ws.Range("L126").Value = 214300680
$ws.Range("M126").Value = -9428
$ws.Range("N126").Value = -214305620
# Row 132
$ws.Range("H132").Value = 3913.2903
$ws.Range("I132").Value = 3654.3704
$ws.Range("K132").Value = 10963.1112
$ws.Range("M132").Value = -8433.111199999999

